# Adding MPA test automation upload file
# Update the asset/sub-asset number sequences on the "Data" sheet.
# Column K (ANLN1) rows whose current value is 60000288 -> 60000306
# Column L (ANLN2) rows whose current value is 229      -> 236
# Column N (PANL1) rows whose current value is 60000289 -> 60000307
# Column O (PANL2) rows whose current value is 230      -> 237

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$kRows = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)
foreach ($r in $kRows) {
    $ws.Range("K$r").Value = 60000306
}

$lRows = @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)
foreach ($r in $lRows) {
    $ws.Range("L$r").Value = 236
}

$nRows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $nRows) {
    $ws.Range("N$r").Value = 60000307
}

$oRows = @(8, 13, 18, 23, 28)
foreach ($r in $oRows) {
    $ws.Range("O$r").Value = 237
}
